# add XS breakpoints to metadata
#
# The "XS" recombinant rows (strain rows 16 & 17) already carry a
# breakpoints range in column G but were missing the corresponding
# clades_regions annotation in column F. This fills in F16/F17 with the
# rich-text breakpoint summary (matching the styling used by every other
# row in the sheet: the clade names in bold/red and bold/green, the
# numeric ranges in between left as un-bolded default text) and nudges the
# G column breakpoints range by one base pair to 9056:10448.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G16 / G17: breakpoints range shifts from 9055:10447 to 9056:10448 ---
$ws.Range("G16").Value = "9056:10448"
$ws.Range("G17").Value = "9056:10448"

# --- F16: new clades_regions rich-text value ---
$clades = "210:10029|Delta/21J,10449:29742|Omicron/BA.1/21K"
$ws.Range("F16").Value = $clades

# "Delta/21J" runs from character 11 through 19 (length 9) -> bold red
$ws.Range("F16").Characters(11, 9).Font.Bold = $true
$ws.Range("F16").Characters(11, 9).Font.Color = 255

# ",10449:29742|" runs from character 20 through 32 (length 13) -> normal
$ws.Range("F16").Characters(20, 13).Font.Size = 11

# "Omicron/BA.1/21K" runs from character 33 through 48 (length 16) -> bold green
$ws.Range("F16").Characters(33, 16).Font.Bold = $true
$ws.Range("F16").Characters(33, 16).Font.Color = 5287936

# --- F17: identical value/formatting to F16 ---
# Copy/paste (instead of re-building the rich text run by run) so the
# workbook reuses a single shared-string entry for both cells, exactly
# like the rest of the sheet does for repeated clades_regions values.
[void]$ws.Range("F16").Copy()
[void]$ws.Range("F17").PasteSpecial()
$excel.CutCopyMode = $false

# --- sheetView: selection moves from F11 to F14, scrolled so column B is
#     the first visible column ---
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("F14").Select()
